$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" year column header (R4), matching the format already
# used by the neighboring "2020" header cell (Q4).
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 2021

# Add the corresponding data point (R5), based on the neighboring "2020"
# data cell (Q5) formatting, then switch its number format to show one
# decimal place without grouping (this mints a new style/numFmt pair).
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 102.20441221981518
$ws.Range("R5").NumberFormat = "0.0"

# Move the active selection as recorded after the edit.
$ws.Range("S9").Select()
